$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from the end of the "Post-conditions: The
#    balance is not increased" paragraph (Bug 1 section).
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ---------------------------------------------------------------------------
# 2) Bug 2 section - append new runs of text after the trailing ": " run in
#    the Description / Pre-conditions / Post-conditions paragraphs.
# ---------------------------------------------------------------------------

# -- Description: plain run, no explicit run formatting.
$pDescription = $d.Paragraphs.Item(9)
$rDescription = $pDescription.Range
$insDescription = $rDescription.Duplicate
$insDescription.End = $rDescription.End - 1
$insDescription.Collapse(0)
$insDescription.InsertAfter("Player cannot reach betting limit: Limit set to 0, but game ends with player still with 5 (dollars) remaining")

# -- Pre-conditions: run with Arial / 10pt formatting.
$pPre = $d.Paragraphs.Item(10)
$rPre = $pPre.Range
$insPre = $rPre.Duplicate
$insPre.End = $rPre.End - 1
$insPre.Collapse(0)
$preText = [char]0x2019
$preStart = $insPre.Start
$preFull = "Player" + $preText + "s remaining balance is set to the same as their bet"
$insPre.InsertAfter($preFull)
$preRange = $d.Range($preStart, $preStart + $preFull.Length)
$preRange.Font.Name = "Arial"
$preRange.Font.Size = 10

# -- Post-conditions: run with Arial / 10pt formatting, plus the _GoBack
#    bookmark re-created immediately after the new text (collapsed, with no
#    characters between bookmarkStart/bookmarkEnd). This runtime's
#    Bookmarks.Add silently snaps a zero-length range to document position 0,
#    so a one-character throwaway marker is bookmarked instead and then
#    deleted, which leaves the (now collapsed) bookmark correctly anchored.
$pPost = $d.Paragraphs.Item(11)
$rPost = $pPost.Range
$insPost = $rPost.Duplicate
$insPost.End = $rPost.End - 1
$insPost.Collapse(0)
$postStart = $insPost.Start
$postText = "The player was not able to bet"
$insPost.InsertAfter($postText + "X")
$postRange = $d.Range($postStart, $postStart + $postText.Length)
$postRange.Font.Name = "Arial"
$postRange.Font.Size = 10

$postEnd = $pPost.Range.End
$markerStart = $postEnd - 2
$markerEnd = $postEnd - 1
$markerRange = $d.Range($markerStart, $markerEnd)
$markerRange.Bookmarks.Add("_GoBack")
$d.Range($markerStart, $markerEnd).Text = ""

# ---------------------------------------------------------------------------
# 3) Footer: collapse the spell-check-wrapped "Jaan" / " " / "Liiband" runs
#    into a single plain run "Jaan Liiband" (no proofErr markers).
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("Jaan Liiband", $true, $false, $false, $false, $false, $true, 1, $false, "Jaan Liiband", 2)
